$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 24.85000000000045
$ws.Range("H2").Value = 0.000884657415986756
$ws.Range("I2").Value = 0.000884657415986756
$ws.Range("L2").Value = 42.78368000611282
$ws.Range("M2").Value = '[15.001040227337327, 70.56631978488831]'
$ws.Range("N2").Value = 0.003318736916741116
$ws.Range("O2").Value = 0.003318736916741116
$ws.Range("P2").Value = 1.968605606753965
$ws.Range("Q2").Value = '[1.2641844311742716, 2.6730267823336575]'
$ws.Range("R2").Value = [double]"1.107176490044992e-06"
$ws.Range("S2").Value = [double]"1.107176490044992e-06"
$ws.Range("T2").Value = 62.56156397856243
$ws.Range("U2").Value = '[47.325800477578255, 77.79732747954661]'
$ws.Range("V2").Value = [double]"1.397588711427034e-10"
$ws.Range("W2").Value = [double]"1.397588711427034e-10"
$ws.Range("X2").Value = 17.06416416416447
$ws.Range("Y2").Value = 14.27817817817844
$ws.Range("Z2").Value = 19.8501501501505

$ws.Range("F3").Value = 24.85000000000045
$ws.Range("H3").Value = 0.02663817990316186
$ws.Range("I3").Value = 0.02663817990316186
$ws.Range("L3").Value = 35.29946307965034
$ws.Range("M3").Value = '[2.429359282759137, 68.16956687654155]'
$ws.Range("N3").Value = 0.03589508871728508
$ws.Range("O3").Value = 0.03589508871728508
$ws.Range("P3").Value = 1.880552959806503
$ws.Range("Q3").Value = '[0.5471843060306547, 3.213921613582352]'
$ws.Range("R3").Value = 0.006739083232343956
$ws.Range("S3").Value = 0.006739083232343956
$ws.Range("T3").Value = 62.73631481755064
$ws.Range("U3").Value = '[44.51036926496292, 80.96226037013835]'
$ws.Range("V3").Value = [double]"1.281749950621247e-08"
$ws.Range("W3").Value = [double]"1.281749950621247e-08"
$ws.Range("X3").Value = 17.41241241241272
$ws.Range("Y3").Value = 12.13893893893915
$ws.Range("Z3").Value = 22.68588588588629

$ws.Range("F4").Value = 24.85000000000045
$ws.Range("H4").Value = 0.001856336920416757
$ws.Range("I4").Value = 0.001856336920416757
$ws.Range("L4").Value = 38.05220351067722
$ws.Range("M4").Value = '[14.44073803188801, 61.66366898946642]'
$ws.Range("N4").Value = 0.002213030275489247
$ws.Range("O4").Value = 0.002213030275489247
$ws.Range("P4").Value = 1.515763422452733
$ws.Range("Q4").Value = '[0.6478159025420398, 2.383710942363426]'
$ws.Range("R4").Value = 0.00100849462680519
$ws.Range("S4").Value = 0.00100849462680519
$ws.Range("T4").Value = 53.12829086512483
$ws.Range("U4").Value = '[38.672152689127614, 67.58442904112205]'
$ws.Range("V4").Value = [double]"2.591338699176049e-09"
$ws.Range("W4").Value = [double]"2.591338699176049e-09"
$ws.Range("X4").Value = 18.85515515515549
$ws.Range("Y4").Value = 15.4224224224227
$ws.Range("Z4").Value = 22.28788788788829

$ws.Range("F5").Value = 24.85000000000045
$ws.Range("H5").Value = [double]"1.20003863512963e-07"
$ws.Range("I5").Value = [double]"1.20003863512963e-07"
$ws.Range("L5").Value = 66.00121532481074
$ws.Range("M5").Value = '[44.09501065586208, 87.9074199937594]'
$ws.Range("N5").Value = [double]"2.472151554577806e-07"
$ws.Range("O5").Value = [double]"2.472151554577806e-07"
$ws.Range("P5").Value = 1.553500271144502
$ws.Range("Q5").Value = '[1.1635528346628865, 1.9434477076261176]'
$ws.Range("R5").Value = [double]"3.181550578545966e-10"
$ws.Range("S5").Value = [double]"3.181550578545966e-10"
$ws.Range("T5").Value = 67.00852154405884
$ws.Range("U5").Value = '[52.97506142619477, 81.0419816619229]'
$ws.Range("V5").Value = [double]"1.74105174721717e-12"
$ws.Range("W5").Value = [double]"1.74105174721717e-12"
$ws.Range("X5").Value = 18.70590590590624
$ws.Range("Y5").Value = 17.16366366366398
$ws.Range("Z5").Value = 20.24814814814851

$ws.Range("F6").Value = 24.85000000000045
$ws.Range("H6").Value = [double]"2.056946801476389e-05"
$ws.Range("I6").Value = [double]"2.056946801476389e-05"
$ws.Range("L6").Value = 56.73026279202769
$ws.Range("M6").Value = '[27.031313727215547, 86.42921185683983]'
$ws.Range("N6").Value = 0.0003738253742899822
$ws.Range("O6").Value = 0.0003738253742899822
$ws.Range("P6").Value = 0.9119738433844251
$ws.Range("Q6").Value = '[0.42139481039142357, 1.4025528763774266]'
$ws.Range("R6").Value = 0.0005118688946175176
$ws.Range("S6").Value = 0.0005118688946175176
$ws.Range("T6").Value = 68.97972579756085
$ws.Range("U6").Value = '[53.659101970674925, 84.30034962444677]'
$ws.Range("V6").Value = [double]"1.014965889112318e-11"
$ws.Range("W6").Value = [double]"1.014965889112318e-11"
$ws.Range("X6").Value = 21.24314314314352
$ws.Range("Y6").Value = 19.30290290290324
$ws.Range("Z6").Value = 23.1833833833838

$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 24.85000000000045
$ws.Range("H7").Value = 0.05045193481900068
$ws.Range("I7").Value = 0.05045193481900068
$ws.Range("L7").Value = 29.4048473712901
$ws.Range("M7").Value = '[-2.453067521628242, 61.26276226420845]'
$ws.Range("N7").Value = 0.06957142502199609
$ws.Range("O7").Value = 0.06957142502199609
$ws.Range("P7").Value = 1.037763339023655
$ws.Range("Q7").Value = '[-0.7232895999255788, 2.7988162779728896]'
$ws.Range("R7").Value = 0.2415021067940069
$ws.Range("S7").Value = 0.2415021067940069
$ws.Range("T7").Value = 56.31525247702592
$ws.Range("U7").Value = '[39.68783827720893, 72.94266667684292]'
$ws.Range("V7").Value = [double]"1.875075028578976e-08"
$ws.Range("W7").Value = [double]"1.875075028578976e-08"
$ws.Range("X7").Value = 20.74564564564602
$ws.Range("Y7").Value = 13.78068068068093
$ws.Range("Z7").Value = 27.71061061061111

$ws.Range("F8").Value = 24.11000000000033
$ws.Range("H8").Value = 0.01840909328525253
$ws.Range("I8").Value = 0.01840909328525253
$ws.Range("L8").Value = 32.30584706696716
$ws.Range("M8").Value = '[4.174101183945339, 60.43759294998898]'
$ws.Range("N8").Value = 0.02535367450428017
$ws.Range("O8").Value = 0.02535367450428017
$ws.Range("P8").Value = 1.176131784226809
$ws.Range("Q8").Value = '[0.05660527303765228, 2.295658295415965]'
$ws.Range("R8").Value = 0.03991616541364729
$ws.Range("S8").Value = 0.03991616541364729
$ws.Range("T8").Value = 59.85920896249031
$ws.Range("U8").Value = '[44.28047582952544, 75.43794209545518]'
$ws.Range("V8").Value = [double]"8.290672592892179e-10"
$ws.Range("W8").Value = [double]"8.290672592892179e-10"
$ws.Range("X8").Value = 19.59691691691719
$ws.Range("Y8").Value = 15.30104104104125
$ws.Range("Z8").Value = 23.89279279279312

$ws.Range("F9").Value = 24.11000000000033
$ws.Range("H9").Value = 0.0001122090766284556
$ws.Range("I9").Value = 0.0001122090766284556
$ws.Range("L9").Value = 54.66402896861776
$ws.Range("M9").Value = '[23.276766172377933, 86.05129176485758]'
$ws.Range("N9").Value = 0.001037502899803311
$ws.Range("O9").Value = 0.001037502899803311
$ws.Range("P9").Value = 0.9748685912040393
$ws.Range("Q9").Value = '[0.40881586082750054, 1.5409213215805782]'
$ws.Range("R9").Value = 0.001163537437558926
$ws.Range("S9").Value = 0.001163537437558926
$ws.Range("T9").Value = 71.90740765503583
$ws.Range("U9").Value = '[55.46860330989432, 88.34621200017733]'
$ws.Range("V9").Value = [double]"2.353472972060899e-11"
$ws.Range("W9").Value = [double]"2.353472972060899e-11"
$ws.Range("X9").Value = 20.36920920920949
$ws.Range("Y9").Value = 18.19713713713739
$ws.Range("Z9").Value = 22.54128128128159

$ws.Range("F10").Value = 24.11000000000033
$ws.Range("H10").Value = [double]"3.216206428535706e-05"
$ws.Range("I10").Value = [double]"3.216206428535706e-05"
$ws.Range("L10").Value = 51.22205935383251
$ws.Range("M10").Value = '[24.392532852049257, 78.05158585561576]'
$ws.Range("N10").Value = 0.0003761626405514029
$ws.Range("O10").Value = 0.0003761626405514029
$ws.Range("P10").Value = 0.6352369529781168
$ws.Range("Q10").Value = '[0.10692107129334794, 1.1635528346628856]'
$ws.Range("R10").Value = 0.01953991407723343
$ws.Range("S10").Value = 0.01953991407723343
$ws.Range("T10").Value = 64.20008474615449
$ws.Range("U10").Value = '[49.98474089530744, 78.41542859700154]'
$ws.Range("V10").Value = [double]"9.272360657064382e-12"
$ws.Range("W10").Value = [double]"9.272360657064382e-12"
$ws.Range("X10").Value = 21.67245245245275
$ws.Range("Y10").Value = 19.64518518518545
$ws.Range("Z10").Value = 23.69971971972004

$ws.Range("F11").Value = 24.11000000000033
$ws.Range("H11").Value = 0.001046857411584501
$ws.Range("I11").Value = 0.001046857411584501
$ws.Range("L11").Value = 44.76946767811643
$ws.Range("M11").Value = '[15.21033180659579, 74.32860354963708]'
$ws.Range("N11").Value = 0.003822163665074152
$ws.Range("O11").Value = 0.003822163665074152
$ws.Range("P11").Value = 1.050342288587578
$ws.Range("Q11").Value = '[0.3836579616996545, 1.7170266154755023]'
$ws.Range("R11").Value = 0.002717772177496469
$ws.Range("S11").Value = 0.002717772177496469
$ws.Range("T11").Value = 60.36549219170555
$ws.Range("U11").Value = '[44.449164274876196, 76.2818201085349]'
$ws.Range("V11").Value = [double]"1.162248963737511e-09"
$ws.Range("W11").Value = [double]"1.162248963737511e-09"
$ws.Range("X11").Value = 20.07959959959987
$ws.Range("Y11").Value = 17.52138138138162
$ws.Range("Z11").Value = 22.63781781781812

$ws.Range("F12").Value = 24.11000000000033
$ws.Range("H12").Value = [double]"6.863678638779902e-07"
$ws.Range("I12").Value = [double]"6.863678638779902e-07"
$ws.Range("L12").Value = 62.63015442334315
$ws.Range("M12").Value = '[34.586499936521676, 90.67380891016462]'
$ws.Range("N12").Value = [double]"4.789168175012293e-05"
$ws.Range("O12").Value = [double]"4.789168175012293e-05"
$ws.Range("P12").Value = 0.8616580451287312
$ws.Range("Q12").Value = '[0.44655270951926873, 1.2767633807381937]'
$ws.Range("R12").Value = 0.0001322894198982905
$ws.Range("S12").Value = 0.0001322894198982905
$ws.Range("T12").Value = 76.43184357108227
$ws.Range("U12").Value = '[62.26308441928472, 90.60060272287981]'
$ws.Range("V12").Value = [double]"3.641531520770513e-14"
$ws.Range("W12").Value = [double]"3.641531520770513e-14"
$ws.Range("X12").Value = 20.80362362362391
$ws.Range("Y12").Value = 19.21077077077104
$ws.Range("Z12").Value = 22.39647647647679
